$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = -7.384100000000001
$ws.Range("B3").Value = 6.421800000000004
$ws.Range("D3").Value = -7.092199999999997
$ws.Range("B4").Value = 8.8498
$ws.Range("D9").Value = -7.198099999999998
$ws.Range("A11").Value = -21.83710000000001
$ws.Range("A12").Value = -21.59890000000001
$ws.Range("B14").Value = 6.615299999999994
$ws.Range("A15").Value = -21.91459999999999
$ws.Range("D15").Value = -8.220899999999995
$ws.Range("D19").Value = -7.373699999999999
$ws.Range("D20").Value = -7.238300000000004
$ws.Range("D25").Value = -7.874400000000005
$ws.Range("B26").Value = 3.959500000000005
$ws.Range("A27").Value = -21.8558
$ws.Range("D27").Value = -9.068499999999997
$ws.Range("A28").Value = -21.73560000000001
$ws.Range("D28").Value = -8.250399999999999
$ws.Range("D30").Value = -7.425500000000002
$ws.Range("A31").Value = -21.64400000000002
$ws.Range("B31").Value = 4.362300000000004
$ws.Range("A32").Value = -21.43490000000001
$ws.Range("D32").Value = -8.963599999999991
$ws.Range("B35").Value = 8.948600000000001
$ws.Range("A36").Value = -20.01769999999999
$ws.Range("B37").Value = 8.657500000000004
$ws.Range("A38").Value = -19.4947
$ws.Range("B39").Value = 9.051399999999999
$ws.Range("B40").Value = 9.490699999999997
$ws.Range("D44").Value = -7.4324
$ws.Range("B45").Value = 5.8781
$ws.Range("A46").Value = -21.5868
$ws.Range("D47").Value = -7.559099999999997
$ws.Range("B52").Value = 5.5959
$ws.Range("A54").Value = -21.6025
$ws.Range("A55").Value = -22.3244
$ws.Range("A56").Value = -22.0789
$ws.Range("B57").Value = 4.766199999999995
$ws.Range("D58").Value = -8.164999999999994
$ws.Range("D62").Value = -8.404499999999992
$ws.Range("A67").Value = -21.49739999999997
$ws.Range("A69").Value = -21.72379999999997
$ws.Range("A72").Value = -21.8968
$ws.Range("A73").Value = -19.72170000000001
$ws.Range("D77").Value = -5.6921
$ws.Range("D78").Value = -7.366200000000005
$ws.Range("B81").Value = 5.391900000000001
$ws.Range("A83").Value = -21.58629999999999
$ws.Range("B83").Value = 6.487900000000001
$ws.Range("D84").Value = -8.780900000000001
$ws.Range("A86").Value = -22.0219
$ws.Range("D89").Value = -6.388299999999996
$ws.Range("A91").Value = -21.47060000000002
$ws.Range("D91").Value = -6.329299999999996
$ws.Range("D92").Value = -6.232599999999998
$ws.Range("A93").Value = -21.34229999999999
$ws.Range("D96").Value = -7.457400000000004
$ws.Range("A99").Value = -20.31429999999999
$ws.Range("B100").Value = 5.920799999999997
$ws.Range("B102").Value = 8.273300000000001
$ws.Range("D102").Value = -8.090500000000002
